# Apply Sprint 3 Backlog - Burndown updates
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Row 3: Filter recipes / Create UI for planned meals page(Desktop) - Janera
$ws.Range("F3").Value = 1

# Row 5: Filter recipes / Implement functionality for view planned meals(Dekstop) - Janera
$ws.Range("F5").Value = 2.5

# Row 9: task now split between Janera and Matthew; Actual time recorded
$ws.Range("D9").Value = 0.5
$ws.Range("E9").Value = "Matthew"
$ws.Range("G9").Value = 0

# Row 11: task re-assigned/split - estimate reduced from 1 to 0.5, remainder to D11, actual time logged
$ws.Range("C11").Value = 0.5
$ws.Range("D11").Value = 0.5
$ws.Range("E11").Value = "Matthew"
$ws.Range("F11").Value = 0.5
$ws.Range("G11").Value = 0

# Row 13: Actual time logged
$ws.Range("F13").Value = 1.5

# Update the selected cell to reflect the author's final cursor position
$ws.Range("F13").Select()

$wb.Save()
